$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H74").Value = 169166.67
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 251250
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 251250
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -253122
$ws.Range("H77").Value = 169166.67
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 251250
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 1256250
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -1265610
$ws.Range("H127").Value = 928.6667
$ws.Range("I127").Value = 928.6667
$ws.Range("K127").Value = 2786.0001
$ws.Range("M127").Value = 2173.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10006
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H61").Value = 2656.1428
$ws.Range("I61").Value = 1848.8334
$ws.Range("K61").Value = 1848.8334
$ws.Range("M61").Value = -1636.8334
$ws.Range("H136").Value = 2656.1428
$ws.Range("I136").Value = 1848.8334
$ws.Range("K136").Value = 5546.5002
$ws.Range("M136").Value = -2996.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 29999.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("N92").Value = -34991.5
$ws.Range("H103").Value = 8111.2
$ws.Range("J103").Value = 8111.2
$ws.Range("L103").Value = 8111.2
$ws.Range("N103").Value = -10455.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H58").Value = 6484.6665
$ws.Range("I58").Value = 7181.6
$ws.Range("K58").Value = 7181.6
$ws.Range("M58").Value = -6978.6
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H93").Value = 28681.2
$ws.Range("I93").Value = 28681.2
$ws.Range("K93").Value = 28681.2
$ws.Range("M93").Value = -26809.2
$ws.Range("H132").Value = 3533.3333
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 3550
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 10650
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -15710
$ws.Range("H134").Value = 5801.4
$ws.Range("I134").Value = 5801.4
$ws.Range("K134").Value = 17404.2
$ws.Range("M134").Value = -14869.2
$ws.Range("H136").Value = 6484.6665
$ws.Range("I136").Value = 7181.6
$ws.Range("K136").Value = 21544.8
$ws.Range("M136").Value = -18994.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 196.1
$ws.Range("J12").Value = 298.33334
$ws.Range("L12").Value = 895.0000200000001
$ws.Range("N12").Value = -1241.00002
$ws.Range("H13").Value = 284.125
$ws.Range("I13").Value = 437.8
$ws.Range("J13").Value = 28
$ws.Range("K13").Value = 1313.4
$ws.Range("L13").Value = 84
$ws.Range("M13").Value = -1145.4
$ws.Range("N13").Value = -420
$ws.Range("H39").Value = 6250
$ws.Range("H56").Value = 12010
$ws.Range("I56").Value = 12010
$ws.Range("K56").Value = 12010
$ws.Range("M56").Value = -11480
$ws.Range("H68").Value = 2001.7142
$ws.Range("J68").Value = 2001.7142
$ws.Range("L68").Value = 6005.142599999999
$ws.Range("N68").Value = -7627.142599999999
$ws.Range("H71").Value = 2001.7142
$ws.Range("J71").Value = 2001.7142
$ws.Range("L71").Value = 18015.4278
$ws.Range("N71").Value = -26127.4278
$ws.Range("H107").Value = 463.5625
$ws.Range("I107").Value = 363.91666
$ws.Range("K107").Value = 1091.74998
$ws.Range("M107").Value = 828.2500199999999
$ws.Range("H122").Value = 703.7778
$ws.Range("I122").Value = 672.75
$ws.Range("J122").Value = 728.6
$ws.Range("K122").Value = 6054.75
$ws.Range("L122").Value = 6557.400000000001
$ws.Range("M122").Value = -3604.75
$ws.Range("N122").Value = -11457.4
$ws.Range("H131").Value = 1609.3334
$ws.Range("J131").Value = 1997.5
$ws.Range("L131").Value = 5992.5
$ws.Range("N131").Value = -16072.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1791.5
$ws.Range("I3").Value = 350
$ws.Range("J3").Value = 3233
$ws.Range("K3").Value = 350
$ws.Range("L3").Value = 3233
$ws.Range("M3").Value = -234
$ws.Range("N3").Value = -3465
$ws.Range("H4").Value = 1562.5
$ws.Range("J4").Value = 416.66666
$ws.Range("L4").Value = 416.66666
$ws.Range("N4").Value = -640.66666
$ws.Range("H74").Value = 38331
$ws.Range("J74").Value = 38331
$ws.Range("L74").Value = 38331
$ws.Range("N74").Value = -40203
$ws.Range("H77").Value = 38331
$ws.Range("J77").Value = 38331
$ws.Range("L77").Value = 114993
$ws.Range("N77").Value = -124353

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2751.7368
$ws.Range("I22").Value = 2596.2144
$ws.Range("J22").Value = 3187.2
$ws.Range("K22").Value = 2596.2144
$ws.Range("L22").Value = 3187.2
$ws.Range("M22").Value = -2301.2144
$ws.Range("N22").Value = -3777.2
$ws.Range("H27").Value = 2751.7368
$ws.Range("I27").Value = 2596.2144
$ws.Range("J27").Value = 3187.2
$ws.Range("K27").Value = 2596.2144
$ws.Range("L27").Value = 3187.2
$ws.Range("M27").Value = -2489.2144
$ws.Range("N27").Value = -3401.2
$ws.Range("H46").Value = 2195.8
$ws.Range("I46").Value = 2142.8572
$ws.Range("J46").Value = 2319.3333
$ws.Range("K46").Value = 2142.8572
$ws.Range("L46").Value = 2319.3333
$ws.Range("M46").Value = -1954.8572
$ws.Range("N46").Value = -2695.3333
$ws.Range("H132").Value = 9918.179
$ws.Range("I132").Value = 8100.5293
$ws.Range("K132").Value = 24301.5879
$ws.Range("M132").Value = -21771.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2039.5
$ws.Range("I126").Value = 2039.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6118.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3648.5
$ws.Range("N126").ClearContents()
